# Refresh the cryptos price/volume table (rows 2-51) plus the two rows whose coin
# identities were swapped (EnergySwap <-> Cronos at rows 48/49), matching the commit:
#   "Updated cryptos list on Wed Sep  6 03:48:42 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text. Values that look like plain numbers are prefixed with a leading
# apostrophe so Excel stores them as literal text (matching the source file's
# inlineStr cells) instead of silently coercing them into floating point numbers,
# e.g. "215.65" -> 215.65000000000001 or "0.06458" -> 6.458E-02.
$updates = @(
    @{ Cell = "D2"; Value = "25.887.49" }
    @{ Cell = "E2"; Value = "  +0.53%  " }
    @{ Cell = "D3"; Value = "1.639.92" }
    @{ Cell = "E3"; Value = "  +1.25%  " }
    @{ Cell = "E4"; Value = "  +0.08%  " }
    @{ Cell = "D5"; Value = "'215.65" }
    @{ Cell = "E5"; Value = "  +0.44%  " }
    @{ Cell = "D6"; Value = "'0.5079" }
    @{ Cell = "E6"; Value = "  -0.07%  " }
    @{ Cell = "E7"; Value = "  +0.07%  " }
    @{ Cell = "D8"; Value = "'0.2598" }
    @{ Cell = "E8"; Value = "  +1.33%  " }
    @{ Cell = "D9"; Value = "'0.06458" }
    @{ Cell = "E9"; Value = "  +1.63%  " }
    @{ Cell = "D10"; Value = "'20.23" }
    @{ Cell = "E10"; Value = "  +5.09%  " }
    @{ Cell = "D11"; Value = "'0.07815" }
    @{ Cell = "E11"; Value = "  +0.59%  " }
    @{ Cell = "D12"; Value = "1.665.67" }
    @{ Cell = "E12"; Value = "  +2.89%  " }
    @{ Cell = "D13"; Value = "'4.268" }
    @{ Cell = "E13"; Value = "  +0.73%  " }
    @{ Cell = "D14"; Value = "1.866.08" }
    @{ Cell = "E14"; Value = "  +1.30%  " }
    @{ Cell = "D15"; Value = "'0.5660" }
    @{ Cell = "E15"; Value = "  +2.18%  " }
    @{ Cell = "D16"; Value = "0.0₅7698" }
    @{ Cell = "E16"; Value = "  +2.60%  " }
    @{ Cell = "D17"; Value = "'63.48" }
    @{ Cell = "E17"; Value = "  +0.00%  " }
    @{ Cell = "D18"; Value = "25.906.32" }
    @{ Cell = "E18"; Value = "  +0.65%  " }
    @{ Cell = "E19"; Value = "  +0.08%  " }
    @{ Cell = "D20"; Value = "'194.64" }
    @{ Cell = "E20"; Value = "  +0.68%  " }
    @{ Cell = "D21"; Value = "'4.395" }
    @{ Cell = "E21"; Value = "  +1.32%  " }
    @{ Cell = "D22"; Value = "'9.987" }
    @{ Cell = "E22"; Value = "  +2.55%  " }
    @{ Cell = "D23"; Value = "'6.245" }
    @{ Cell = "E23"; Value = "  +4.80%  " }
    @{ Cell = "D24"; Value = "'1.004" }
    @{ Cell = "E24"; Value = "  +0.10%  " }
    @{ Cell = "D25"; Value = "'1.757" }
    @{ Cell = "E25"; Value = "  -4.05%  " }
    @{ Cell = "D26"; Value = "'138.57" }
    @{ Cell = "E26"; Value = "  -1.52%  " }
    @{ Cell = "D27"; Value = "'0.1228" }
    @{ Cell = "E27"; Value = "  -2.59%  " }
    @{ Cell = "D28"; Value = "'6.859" }
    @{ Cell = "E28"; Value = "  +2.15%  " }
    @{ Cell = "D29"; Value = "'15.56" }
    @{ Cell = "E29"; Value = "  +1.11%  " }
    @{ Cell = "E30"; Value = "  +0.93%  " }
    @{ Cell = "D31"; Value = "'0.04979" }
    @{ Cell = "E31"; Value = "  +2.82%  " }
    @{ Cell = "D32"; Value = "'3.323" }
    @{ Cell = "E32"; Value = "  +1.07%  " }
    @{ Cell = "D33"; Value = "'3.259" }
    @{ Cell = "E33"; Value = "  +2.84%  " }
    @{ Cell = "E34"; Value = "  +2.04%  " }
    @{ Cell = "D35"; Value = "'2.393" }
    @{ Cell = "E35"; Value = "  +1.19%  " }
    @{ Cell = "D36"; Value = "'0.9095" }
    @{ Cell = "D37"; Value = "'2.579" }
    @{ Cell = "E37"; Value = "  +1.72%  " }
    @{ Cell = "D38"; Value = "'0.5532" }
    @{ Cell = "E38"; Value = "  +1.36%  " }
    @{ Cell = "D39"; Value = "1.127.56" }
    @{ Cell = "E39"; Value = "  +0.41%  " }
    @{ Cell = "D40"; Value = "'0.01574" }
    @{ Cell = "E40"; Value = "  +1.17%  " }
    @{ Cell = "D41"; Value = "'0.9994" }
    @{ Cell = "E41"; Value = "  -0.77%  " }
    @{ Cell = "E42"; Value = "  -1.01%  " }
    @{ Cell = "D43"; Value = "'99.60" }
    @{ Cell = "E43"; Value = "  +2.66%  " }
    @{ Cell = "D44"; Value = "'0.8012" }
    @{ Cell = "E44"; Value = "  +1.26%  " }
    @{ Cell = "E45"; Value = "  -3.98%  " }
    @{ Cell = "D46"; Value = "'55.68" }
    @{ Cell = "E46"; Value = "  +2.08%  " }
    @{ Cell = "D47"; Value = "'0.4238" }
    @{ Cell = "E47"; Value = "  -4.10%  " }
    @{ Cell = "B48"; Value = "Cronos" }
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro" }
    @{ Cell = "D48"; Value = "'0.05043" }
    @{ Cell = "E48"; Value = "  -0.46%  " }
    @{ Cell = "B49"; Value = "EnergySwap" }
    @{ Cell = "C49"; Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens" }
    @{ Cell = "D49"; Value = "'7.646" }
    @{ Cell = "E49"; Value = "  +1.37%  " }
    @{ Cell = "D50"; Value = "'0.9998" }
    @{ Cell = "E50"; Value = "  +0.05%  " }
    @{ Cell = "E51"; Value = "  +0.09%  " }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

Write-Host "Applied $($updates.Count) cell updates"
